$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.455.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.249.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.631'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.89'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.31%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.623'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.75'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0949'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.590.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.62'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.857'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.245.69'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.248.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  +2.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +33.89%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.61'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.121'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.95'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.73%  '
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.51'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0318'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.15%  '
$ws.Range("E40").Value = '  -2.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '63.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.55%  '
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '108.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.102'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.998'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.19'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.55%  '
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.13'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.23%  '
